$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Fecha): weekly shift - each record date moves down one slot, newest record added at top ---
$ws.Range("D82").Value = 44582
$ws.Range("D83").Value = 44582
$ws.Range("D84").Value = 44160
$ws.Range("D85").Value = 44160
$ws.Range("D86").Value = 44308
$ws.Range("D87").Value = 44308
$ws.Range("D88").Value = 44320
$ws.Range("D89").Value = 44320
$ws.Range("D90").Value = 44306
$ws.Range("D91").Value = 44306
$ws.Range("D92").Value = 44316
$ws.Range("D93").Value = 44316
$ws.Range("D94").Value = 44460
$ws.Range("D95").Value = 44460
$ws.Range("D96").Value = 44272
$ws.Range("D97").Value = 44272
$ws.Range("D98").Value = 44313
$ws.Range("D99").Value = 44313
$ws.Range("D100").Value = 44334
$ws.Range("D101").Value = 44334
$ws.Range("D102").Value = 44405
$ws.Range("D103").Value = 44405
$ws.Range("D104").Value = 44280
$ws.Range("D105").Value = 44280
$ws.Range("D106").Value = 44330
$ws.Range("D107").Value = 44330
$ws.Range("D108").Value = 44239
$ws.Range("D109").Value = 44239
$ws.Range("D110").Value = 44476
$ws.Range("D111").Value = 44476
$ws.Range("D112").Value = 44250
$ws.Range("D113").Value = 44250
$ws.Range("D114").Value = 44488
$ws.Range("D115").Value = 44488
$ws.Range("D116").Value = 44341
$ws.Range("D117").Value = 44341
$ws.Range("D118").Value = 44278
$ws.Range("D119").Value = 44278
$ws.Range("D120").Value = 44194
$ws.Range("D121").Value = 44194
$ws.Range("D122").Value = 44490
$ws.Range("D123").Value = 44490
$ws.Range("D124").Value = 44525
$ws.Range("D125").Value = 44525
$ws.Range("D126").Value = 44327
$ws.Range("D127").Value = 44327

# --- Column O (Origen): region values follow the same shift for rows 110-113 ---
$ws.Range("O110").Value = "Región de Ñuble"
$ws.Range("O111").Value = "Región de Ñuble"
$ws.Range("O112").Value = "Región de Arica y Parinacota"
$ws.Range("O113").Value = "Región de Arica y Parinacota"

# --- Column J (Volumen): values follow the same shift for rows 116-119 ---
$ws.Range("J116").Value = 200
$ws.Range("J117").Value = 100
$ws.Range("J118").Value = 300
$ws.Range("J119").Value = 150

# --- New rows 128 and 129 (the record pushed out of row 126/127 by the shift) ---
# Row 128
$ws.Range("A128").Value = 11
$ws.Range("B128").Value = "Vega Monumental Concepción"
$ws.Range("C128").Value = "Bíobío"
$ws.Range("D128").Value = 44512
$ws.Range("E128").Value = 8
$ws.Range("F128").Value = 100112044
$ws.Range("G128").Value = "Perejil"
$ws.Range("H128").Value = "Sin especificar"
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 200
$ws.Range("K128").Value = 600
$ws.Range("L128").Value = 700
$ws.Range("M128").Value = 650
$ws.Range("N128").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O128").Value = "Región de Ñuble"
$ws.Range("P128").Value = 650
$ws.Range("Q128").Value = 1
$ws.Range("R128").Value = "Hortaliza"
$ws.Range("D128").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 129
$ws.Range("A129").Value = 11
$ws.Range("B129").Value = "Vega Monumental Concepción"
$ws.Range("C129").Value = "Bíobío"
$ws.Range("D129").Value = 44512
$ws.Range("E129").Value = 8
$ws.Range("F129").Value = 100112044
$ws.Range("G129").Value = "Perejil"
$ws.Range("H129").Value = "Sin especificar"
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 100
$ws.Range("K129").Value = 500
$ws.Range("L129").Value = 500
$ws.Range("M129").Value = 500
$ws.Range("N129").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O129").Value = "Región de Ñuble"
$ws.Range("P129").Value = 500
$ws.Range("Q129").Value = 1
$ws.Range("R129").Value = "Hortaliza"
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
